$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.458056666666667
$ws.Range("H2").Value = 4.374169999999999
$ws.Range("I2").Value = 0.2323568509805328
$ws.Range("J2").Value = 0.2323568509805327
$ws.Range("M2").Value = 60.538204
$ws.Range("N2").Value = 181.614612
$ws.Range("O2").Value = 0.6123615450168176
$ws.Range("P2").Value = 0.6123615450168176
$ws.Range("Q2").Value = 88.26813193022666
$ws.Range("R2").Value = 794.4131873720398
$ws.Range("S2").Value = 0.1422864002616815
$ws.Range("T2").Value = 0.1422864002616815
$ws.Range("G3").Value = 1.458056666666667
$ws.Range("H3").Value = 4.374169999999999
$ws.Range("I3").Value = 0.2323568509805328
$ws.Range("J3").Value = 0.2323568509805327
$ws.Range("O3").Value = 0.1096681363892149
$ws.Range("P3").Value = 0.1096681363892149
$ws.Range("Q3").Value = 15.80798404164889
$ws.Range("R3").Value = 142.27185637484
$ws.Range("S3").Value = 0.02548214282430154
$ws.Range("T3").Value = 0.02548214282430154
$ws.Range("G4").Value = 1.458056666666667
$ws.Range("H4").Value = 4.374169999999999
$ws.Range("I4").Value = 0.2323568509805328
$ws.Range("J4").Value = 0.2323568509805327
$ws.Range("M4").Value = 8.850437666666666
$ws.Range("N4").Value = 26.551313
$ws.Range("O4").Value = 0.08952475173586316
$ws.Range("P4").Value = 0.08952475173586316
$ws.Range("Q4").Value = 12.90443964280111
$ws.Range("R4").Value = 116.13995678521
$ws.Range("S4").Value = 0.02080168939815915
$ws.Range("T4").Value = 0.02080168939815915
$ws.Range("G5").Value = 1.458056666666667
$ws.Range("H5").Value = 4.374169999999999
$ws.Range("I5").Value = 0.2323568509805328
$ws.Range("J5").Value = 0.2323568509805327
$ws.Range("M5").Value = 1.757142
$ws.Range("N5").Value = 5.271426
$ws.Range("O5").Value = 0.01777400251143792
$ws.Range("P5").Value = 0.01777400251143792
$ws.Range("Q5").Value = 2.56201260738
$ws.Range("R5").Value = 23.05811346642
$ws.Range("S5").Value = 0.004129911252877795
$ws.Range("T5").Value = 0.004129911252877795
$ws.Range("G6").Value = 1.458056666666667
$ws.Range("H6").Value = 4.374169999999999
$ws.Range("I6").Value = 0.2323568509805328
$ws.Range("J6").Value = 0.2323568509805327
$ws.Range("M6").Value = 16.87263033333333
$ws.Range("N6").Value = 50.617891
$ws.Range("O6").Value = 0.1706715643466665
$ws.Range("P6").Value = 0.1706715643466665
$ws.Range("Q6").Value = 24.60125114171889
$ws.Range("R6").Value = 221.41126027547
$ws.Range("S6").Value = 0.03965670724351279
$ws.Range("T6").Value = 0.03965670724351279
$ws.Range("H7").Value = 5.708772
$ws.Range("I7").Value = 0.3032511962008422
$ws.Range("J7").Value = 0.3032511962008422
$ws.Range("M7").Value = 60.538204
$ws.Range("N7").Value = 181.614612
$ws.Range("O7").Value = 0.6123615450168176
$ws.Range("P7").Value = 0.6123615450168176
$ws.Range("Q7").Value = 115.199601308496
$ws.Range("R7").Value = 1036.796411776464
$ws.Range("S7").Value = 0.1856993710337458
$ws.Range("T7").Value = 0.1856993710337458
$ws.Range("H8").Value = 5.708772
$ws.Range("I8").Value = 0.3032511962008422
$ws.Range("J8").Value = 0.3032511962008422
$ws.Range("O8").Value = 0.1096681363892149
$ws.Range("P8").Value = 0.1096681363892149
$ws.Range("S8").Value = 0.03325699354514652
$ws.Range("T8").Value = 0.03325699354514653
$ws.Range("H9").Value = 5.708772
$ws.Range("I9").Value = 0.3032511962008422
$ws.Range("J9").Value = 0.3032511962008422
$ws.Range("M9").Value = 8.850437666666666
$ws.Range("N9").Value = 26.551313
$ws.Range("O9").Value = 0.08952475173586316
$ws.Range("P9").Value = 0.08952475173586316
$ws.Range("Q9").Value = 16.841710246404
$ws.Range("R9").Value = 151.575392217636
$ws.Range("S9").Value = 0.02714848805348393
$ws.Range("T9").Value = 0.02714848805348393
$ws.Range("H10").Value = 5.708772
$ws.Range("I10").Value = 0.3032511962008422
$ws.Range("J10").Value = 0.3032511962008422
$ws.Range("M10").Value = 1.757142
$ws.Range("N10").Value = 5.271426
$ws.Range("O10").Value = 0.01777400251143792
$ws.Range("P10").Value = 0.01777400251143792
$ws.Range("Q10").Value = 3.343707683208
$ws.Range("R10").Value = 30.093369148872
$ws.Range("S10").Value = 0.005389987522870322
$ws.Range("T10").Value = 0.005389987522870322
$ws.Range("H11").Value = 5.708772
$ws.Range("I11").Value = 0.3032511962008422
$ws.Range("J11").Value = 0.3032511962008422
$ws.Range("M11").Value = 16.87263033333333
$ws.Range("N11").Value = 50.617891
$ws.Range("O11").Value = 0.1706715643466665
$ws.Range("P11").Value = 0.1706715643466665
$ws.Range("Q11").Value = 32.107333204428
$ws.Range("R11").Value = 288.965998839852
$ws.Range("S11").Value = 0.05175635604559562
$ws.Range("T11").Value = 0.05175635604559562
$ws.Range("G12").Value = 2.914094333333333
$ws.Range("H12").Value = 8.742283
$ws.Range("I12").Value = 0.4643919528186251
$ws.Range("J12").Value = 0.4643919528186251
$ws.Range("M12").Value = 60.538204
$ws.Range("N12").Value = 181.614612
$ws.Range("O12").Value = 0.6123615450168176
$ws.Range("P12").Value = 0.6123615450168176
$ws.Range("Q12").Value = 176.4140372265773
$ws.Range("R12").Value = 1587.726335039196
$ws.Range("S12").Value = 0.2843757737213903
$ws.Range("T12").Value = 0.2843757737213903
$ws.Range("G13").Value = 2.914094333333333
$ws.Range("H13").Value = 8.742283
$ws.Range("I13").Value = 0.4643919528186251
$ws.Range("J13").Value = 0.4643919528186251
$ws.Range("O13").Value = 0.1096681363892149
$ws.Range("P13").Value = 0.1096681363892149
$ws.Range("Q13").Value = 31.59407845410178
$ws.Range("R13").Value = 284.346706086916
$ws.Range("S13").Value = 0.05092900001976681
$ws.Range("T13").Value = 0.05092900001976682
$ws.Range("G14").Value = 2.914094333333333
$ws.Range("H14").Value = 8.742283
$ws.Range("I14").Value = 0.4643919528186251
$ws.Range("J14").Value = 0.4643919528186251
$ws.Range("M14").Value = 8.850437666666666
$ws.Range("N14").Value = 26.551313
$ws.Range("O14").Value = 0.08952475173586316
$ws.Range("P14").Value = 0.08952475173586316
$ws.Range("Q14").Value = 25.79101025195322
$ws.Range("R14").Value = 232.119092267579
$ws.Range("S14").Value = 0.04157457428422009
$ws.Range("T14").Value = 0.04157457428422009
$ws.Range("G15").Value = 2.914094333333333
$ws.Range("H15").Value = 8.742283
$ws.Range("I15").Value = 0.4643919528186251
$ws.Range("J15").Value = 0.4643919528186251
$ws.Range("M15").Value = 1.757142
$ws.Range("N15").Value = 5.271426
$ws.Range("O15").Value = 0.01777400251143792
$ws.Range("P15").Value = 0.01777400251143792
$ws.Range("Q15").Value = 5.120477545062
$ws.Range("R15").Value = 46.084297905558
$ws.Range("S15").Value = 0.008254103735689799
$ws.Range("T15").Value = 0.008254103735689799
$ws.Range("G16").Value = 2.914094333333333
$ws.Range("H16").Value = 8.742283
$ws.Range("I16").Value = 0.4643919528186251
$ws.Range("J16").Value = 0.4643919528186251
$ws.Range("M16").Value = 16.87263033333333
$ws.Range("N16").Value = 50.617891
$ws.Range("O16").Value = 0.1706715643466665
$ws.Range("P16").Value = 0.1706715643466665
$ws.Range("Q16").Value = 49.16843644279478
$ws.Range("R16").Value = 442.515927985153
$ws.Range("S16").Value = 0.07925850105755806
